$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update StartingInventories (C) and SetupCosts (E) ---
$wsProd = $wb.Worksheets.Item("Productdata")

$wsProd.Range("C2").Value = 0
$wsProd.Range("E2").Value = 0.104

$wsProd.Range("C3").Value = 0
$wsProd.Range("E3").Value = 0.09920000000000001

$wsProd.Range("C4").Value = 0
$wsProd.Range("E4").Value = 0.0384

$wsProd.Range("C5").Value = 0
$wsProd.Range("E5").Value = 0.016

$wsProd.Range("C6").Value = 0
$wsProd.Range("E6").Value = 0.0288

$wsProd.Range("C7").Value = 5
$wsProd.Range("E7").Value = 0.052

$wsProd.Range("C8").Value = 5
$wsProd.Range("E8").Value = 0.1016

$wsProd.Range("C9").Value = 5
$wsProd.Range("E9").Value = 0.04960000000000001

# --- Capacity sheet: update capacity values (B) ---
$wsCap = $wb.Worksheets.Item("Capacity")

$wsCap.Range("B2").Value = 100
$wsCap.Range("B3").Value = 40
$wsCap.Range("B4").Value = 160
$wsCap.Range("B5").Value = 200
$wsCap.Range("B6").Value = 160
$wsCap.Range("B7").Value = 10
$wsCap.Range("B8").Value = 10
$wsCap.Range("B9").Value = 30

# --- ProcessingTime sheet: update processing times ---
$wsProc = $wb.Worksheets.Item("ProcessingTime")

$wsProc.Range("C3").Value = 2
$wsProc.Range("F6").Value = 4
$wsProc.Range("I9").Value = 3

$wb.Save()
